$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pantry1")

$ws.Range("A13").Value = "pudding mix"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "packets"

$ws.Range("A14").Select()
